# Source data extraction from sources 1 through 9
# Fill in placeholder "N/A" values for the data columns (TPS, Energy Use per
# Transaction, Nakamoto Coefficient, % of nodes required to take over network,
# Strengths, Weaknesses) for every consensus mechanism row (2-10) on the
# "Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Columns B:G, rows 2-10 -> "N/A"
$ws.Range("B2:G10").Value = "N/A"

# Column E previously carried the "Percent" cell style (reserved for a
# numeric percentage that was never filled in); now that it holds the same
# "N/A" text as the rest of the row, drop back to the default style.
$ws.Range("E2:E10").Style = "Normal"

# The "Percent" cell style is no longer referenced anywhere in the workbook,
# so remove its definition entirely.
$wb.Styles.Item("Percent").Delete()

# Restore the active cell/selection recorded in the sheet at save time.
$ws.Range("H25").Select()
